$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 384.8421
$ws.Cells.Item(33, 9).Value = 395
$ws.Cells.Item(33, 10).Value = 298.5
$ws.Cells.Item(33, 11).Value = 395
$ws.Cells.Item(33, 12).Value = 298.5
$ws.Cells.Item(33, 13).Value = -166
$ws.Cells.Item(33, 14).Value = -756.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 6584.9414
$ws.Cells.Item(40, 9).Value = 4887.6665
$ws.Cells.Item(40, 11).Value = 4887.6665
$ws.Cells.Item(40, 13).Value = -4712.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(42, 8).Value = 547.3333
$ws.Cells.Item(42, 9).Value = 583.5
$ws.Cells.Item(42, 11).Value = 1750.5
$ws.Cells.Item(42, 13).Value = -1520.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 1364.1818
$ws.Cells.Item(80, 9).Value = 715.3333
$ws.Cells.Item(80, 10).Value = 1607.5
$ws.Cells.Item(80, 11).Value = 2145.9999
$ws.Cells.Item(80, 12).Value = 4822.5
$ws.Cells.Item(80, 13).Value = -1147.9999
$ws.Cells.Item(80, 14).Value = -6818.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(83, 8).Value = 1364.1818
$ws.Cells.Item(83, 9).Value = 715.3333
$ws.Cells.Item(83, 10).Value = 1607.5
$ws.Cells.Item(83, 11).Value = 6437.9997
$ws.Cells.Item(83, 12).Value = 14467.5
$ws.Cells.Item(83, 13).Value = -1445.9997
$ws.Cells.Item(83, 14).Value = -24451.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 5212.1333
$ws.Cells.Item(98, 9).Value = 5405.9287
$ws.Cells.Item(98, 11).Value = 5405.9287
$ws.Cells.Item(98, 13).Value = -3907.9287

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 5212.1333
$ws.Cells.Item(122, 9).Value = 5405.9287
$ws.Cells.Item(122, 11).Value = 16217.7861
$ws.Cells.Item(122, 13).Value = -13767.7861

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 3145.5386
$ws.Cells.Item(132, 9).Value = 3145.5386
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 9436.6158
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = $null
$ws.Cells.Item(132, 14).Value = -6906.6158

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 1001.3333
$ws.Cells.Item(135, 9).Value = 873.75
$ws.Cells.Item(135, 11).Value = 7863.75
$ws.Cells.Item(135, 13).Value = -5328.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 355.125
$ws.Cells.Item(5, 9).Value = 316.33334
$ws.Cells.Item(5, 10).Value = 471.5
$ws.Cells.Item(5, 11).Value = 316.33334
$ws.Cells.Item(5, 12).Value = 471.5
$ws.Cells.Item(5, 13).Value = -204.33334
$ws.Cells.Item(5, 14).Value = -695.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1533.9615
$ws.Cells.Item(32, 9).Value = 1551.56
$ws.Cells.Item(32, 11).Value = 1551.56
$ws.Cells.Item(32, 13).Value = -1264.56

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 739.25
$ws.Cells.Item(110, 9).Value = 653.1667
$ws.Cells.Item(110, 10).Value = 997.5
$ws.Cells.Item(110, 11).Value = 653.1667
$ws.Cells.Item(110, 12).Value = 997.5
$ws.Cells.Item(110, 13).Value = 1391.8333
$ws.Cells.Item(110, 14).Value = -5087.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 3136.3333
$ws.Cells.Item(122, 9).Value = 2791.2
$ws.Cells.Item(122, 11).Value = 8373.599999999999
$ws.Cells.Item(122, 13).Value = -5923.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(123, 8).Value = 49999
$ws.Cells.Item(123, 10).Value = 49999
$ws.Cells.Item(123, 12).Value = 49999
$ws.Cells.Item(123, 14).Value = -59799

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2973.7
$ws.Cells.Item(132, 9).Value = 3387.4443
$ws.Cells.Item(132, 10).Value = 2635.182
$ws.Cells.Item(132, 11).Value = 10162.3329
$ws.Cells.Item(132, 12).Value = 7905.545999999999
$ws.Cells.Item(132, 13).Value = -7632.332900000001
$ws.Cells.Item(132, 14).Value = -12965.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 355.125
$ws.Cells.Item(4, 9).Value = 316.33334
$ws.Cells.Item(4, 10).Value = 471.5
$ws.Cells.Item(4, 11).Value = 316.33334
$ws.Cells.Item(4, 12).Value = 471.5
$ws.Cells.Item(4, 13).Value = -201.33334
$ws.Cells.Item(4, 14).Value = -701.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 52002900
$ws.Cells.Item(105, 10).Value = 83336830
$ws.Cells.Item(105, 12).Value = 83336830
$ws.Cells.Item(105, 14).Value = -83340324

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 4274943
$ws.Cells.Item(107, 9).Value = 5129528.5
$ws.Cells.Item(107, 10).Value = 2017.3334
$ws.Cells.Item(107, 11).Value = 5129528.5
$ws.Cells.Item(107, 12).Value = 2017.3334
$ws.Cells.Item(107, 13).Value = -5127608.5
$ws.Cells.Item(107, 14).Value = -5857.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 11985.117
$ws.Cells.Item(134, 9).Value = 3573.6
$ws.Cells.Item(134, 10).Value = 24001.572
$ws.Cells.Item(134, 11).Value = 10720.8
$ws.Cells.Item(134, 12).Value = 72004.716
$ws.Cells.Item(134, 13).Value = -8185.799999999999
$ws.Cells.Item(134, 14).Value = -77074.716

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 448.69232
$ws.Cells.Item(7, 9).Value = 410.625
$ws.Cells.Item(7, 10).Value = 509.6
$ws.Cells.Item(7, 11).Value = 410.625
$ws.Cells.Item(7, 12).Value = 509.6
$ws.Cells.Item(7, 13).Value = -297.625
$ws.Cells.Item(7, 14).Value = -735.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4104.553
$ws.Cells.Item(31, 9).Value = 4854.769
$ws.Cells.Item(31, 10).Value = 3817.7058
$ws.Cells.Item(31, 11).Value = 4854.769
$ws.Cells.Item(31, 12).Value = 3817.7058
$ws.Cells.Item(31, 13).Value = -4559.769
$ws.Cells.Item(31, 14).Value = -4407.7058

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 4104.553
$ws.Cells.Item(34, 9).Value = 4854.769
$ws.Cells.Item(34, 10).Value = 3817.7058
$ws.Cells.Item(34, 11).Value = 4854.769
$ws.Cells.Item(34, 12).Value = 3817.7058
$ws.Cells.Item(34, 13).Value = -4652.769
$ws.Cells.Item(34, 14).Value = -4221.7058

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1829.0488
$ws.Cells.Item(58, 9).Value = 1257.125
$ws.Cells.Item(58, 10).Value = 2636.4707
$ws.Cells.Item(58, 11).Value = 1257.125
$ws.Cells.Item(58, 12).Value = 2636.4707
$ws.Cells.Item(58, 13).Value = -1054.125
$ws.Cells.Item(58, 14).Value = -3042.4707

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 10003249
$ws.Cells.Item(62, 9).Value = 12503373
$ws.Cells.Item(62, 10).Value = 2753
$ws.Cells.Item(62, 11).Value = 12503373
$ws.Cells.Item(62, 12).Value = 2753
$ws.Cells.Item(62, 13).Value = -12502749
$ws.Cells.Item(62, 14).Value = -4001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 10003249
$ws.Cells.Item(65, 9).Value = 12503373
$ws.Cells.Item(65, 10).Value = 2753
$ws.Cells.Item(65, 11).Value = 62516865
$ws.Cells.Item(65, 12).Value = 13765
$ws.Cells.Item(65, 13).Value = -62513745
$ws.Cells.Item(65, 14).Value = -20005

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 1829.0488
$ws.Cells.Item(136, 9).Value = 1257.125
$ws.Cells.Item(136, 10).Value = 2636.4707
$ws.Cells.Item(136, 11).Value = 3771.375
$ws.Cells.Item(136, 12).Value = 7909.4121
$ws.Cells.Item(136, 13).Value = -1221.375
$ws.Cells.Item(136, 14).Value = -13009.4121

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 416.33334
$ws.Cells.Item(34, 9).Value = 416.33334
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 1249.00002
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = $null
$ws.Cells.Item(34, 14).Value = -1165.00002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 5979.2
$ws.Cells.Item(39, 9).Value = 500
$ws.Cells.Item(39, 11).Value = 1500
$ws.Cells.Item(39, 13).Value = -1206

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 5499
$ws.Cells.Item(55, 9).Value = 2999.5
$ws.Cells.Item(55, 11).Value = 8998.5
$ws.Cells.Item(55, 13).Value = -8821.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 19589.455
$ws.Cells.Item(80, 9).Value = 1166.6666
$ws.Cells.Item(80, 10).Value = 26498
$ws.Cells.Item(80, 11).Value = 3499.9998
$ws.Cells.Item(80, 12).Value = 79494
$ws.Cells.Item(80, 13).Value = -2563.9998
$ws.Cells.Item(80, 14).Value = -81366

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83, 8).Value = 19589.455
$ws.Cells.Item(83, 9).Value = 1166.6666
$ws.Cells.Item(83, 10).Value = 26498
$ws.Cells.Item(83, 11).Value = 10499.9994
$ws.Cells.Item(83, 12).Value = 238482
$ws.Cells.Item(83, 13).Value = -5819.999400000001
$ws.Cells.Item(83, 14).Value = -247842

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 1652.0454
$ws.Cells.Item(122, 9).Value = 304
$ws.Cells.Item(122, 11).Value = 2736
$ws.Cells.Item(122, 13).Value = -286

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 5000
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 12).Value = 45000
$ws.Cells.Item(132, 14).Value = -50060

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 20005
$ws.Cells.Item(19, 9).Value = 20005
$ws.Cells.Item(19, 11).Value = 20005
$ws.Cells.Item(19, 13).Value = -19717

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5273.3687
$ws.Cells.Item(70, 9).Value = 4599.7144
$ws.Cells.Item(70, 11).Value = 4599.7144
$ws.Cells.Item(70, 13).Value = -4329.7144

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 5273.3687
$ws.Cells.Item(73, 9).Value = 4599.7144
$ws.Cells.Item(73, 11).Value = 4599.7144
$ws.Cells.Item(73, 13).Value = -3663.7144

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2484.7
$ws.Cells.Item(102, 9).Value = 1978.2
$ws.Cells.Item(102, 10).Value = 4004.2
$ws.Cells.Item(102, 11).Value = 1978.2
$ws.Cells.Item(102, 12).Value = 4004.2
$ws.Cells.Item(102, 13).Value = -356.2
$ws.Cells.Item(102, 14).Value = -7248.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3673.3333
$ws.Cells.Item(122, 9).Value = 2966.6428
$ws.Cells.Item(122, 11).Value = 8899.928400000001
$ws.Cells.Item(122, 13).Value = -6449.928400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2025.5333
$ws.Cells.Item(132, 9).Value = 2025.5333
$ws.Cells.Item(132, 11).Value = 6076.5999
$ws.Cells.Item(132, 13).Value = -3546.5999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1899.8889
$ws.Cells.Item(46, 9).Value = 1449.875
$ws.Cells.Item(46, 11).Value = 1449.875
$ws.Cells.Item(46, 13).Value = -1261.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1275.875
$ws.Cells.Item(55, 9).Value = 1117.8182
$ws.Cells.Item(55, 10).Value = 1623.6
$ws.Cells.Item(55, 11).Value = 1117.8182
$ws.Cells.Item(55, 12).Value = 1623.6
$ws.Cells.Item(55, 13).Value = -944.8181999999999
$ws.Cells.Item(55, 14).Value = -1969.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 12099.846
$ws.Cells.Item(122, 9).Value = 7717.3335
$ws.Cells.Item(122, 11).Value = 23152.0005
$ws.Cells.Item(122, 13).Value = -20702.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4707
$ws.Cells.Item(132, 9).Value = 4144.857
$ws.Cells.Item(132, 10).Value = 5581.4443
$ws.Cells.Item(132, 11).Value = 12434.571
$ws.Cells.Item(132, 12).Value = 16744.3329
$ws.Cells.Item(132, 13).Value = -9904.571
$ws.Cells.Item(132, 14).Value = -21804.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 14706970
$ws.Cells.Item(122, 9).Value = 1155.8125
$ws.Cells.Item(122, 11).Value = 3467.4375
$ws.Cells.Item(122, 13).Value = -1017.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4652.125
$ws.Cells.Item(132, 9).Value = 4302.364
$ws.Cells.Item(132, 10).Value = 8499.5
$ws.Cells.Item(132, 11).Value = 12907.092
$ws.Cells.Item(132, 12).Value = 25498.5
$ws.Cells.Item(132, 13).Value = -10377.092
$ws.Cells.Item(132, 14).Value = -30558.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 216129.58
$ws.Cells.Item(136, 9).Value = 3678.761
$ws.Cells.Item(136, 10).Value = 914182.3
$ws.Cells.Item(136, 11).Value = 11036.283
$ws.Cells.Item(136, 12).Value = 2742546.9
$ws.Cells.Item(136, 13).Value = -8486.282999999999
$ws.Cells.Item(136, 14).Value = -2747646.9
